# Bug fix: remove the stray/duplicate Sap record (SapID 60776e12ccab402de07f4c82)
# that was accidentally entered on the "Saps" sheet. Deleting the entire row
# shifts all subsequent sap records up by one row and keeps the sheet's
# shared-string table tidy (mongoDB backup/bug fix commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Saps")

# Row 2 holds the sap entry with SapID "60776e12ccab402de07f4c82"
$ws.Range("A2").EntireRow.Delete()
